$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("R40" rule), column B: its shared-string value is replaced with
# the literal text "1". Because "1" is numeric-looking, a plain
# `Range.Value = "1"` assignment gets auto-coerced to a *number*, and
# forcing text via NumberFormat/apostrophe directly on the cell would
# permanently swap its cell style (it currently uses style index 23,
# shared with other bordered cells in the table).
#
# To keep the cell's original style intact while still storing "1" as
# text, stage the text value on a scratch cell (format it as Text there),
# copy *values only* onto B11 (this carries the string type without
# touching B11's own formatting/style), then remove all traces of the
# scratch cell.
$dest = $ws.Range("B11")
$scratch = $ws.Range("Z1")

$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$dest.PasteSpecial(-4163)   # xlPasteValues - copies the text value, not the format

$scratch.Clear()            # remove the scratch cell's value + formatting entirely
$excel.CutCopyMode = $false
